# Generate Report for Handoff
# Rotates the localization-status report onto a new source GUID / xliff
# hash and refreshes the handoff/handback timestamps, mirroring a fresh
# CI run of the OpenLocalization handback report.

$wb = $excel.ActiveWorkbook

$oldGuid = "725eaac0-bc1b-4466-a756-334e7858836e"
$newGuid = "4b4c19a2-07bd-4260-a3ed-2607b79109fd"
$oldHash = "7f8cccaf704440dab46dea511e5918db13c18ca2"
$newHash = "1abd6df45ed39ce7c1f45640f125c11961bf1118"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04a8cd3b0e6dc0d77bfbed902929b347a5ef2ee8/e2e/$oldGuid.md"
$hyperlinkColor = 15570276  # RGB(0x64, 0x95, 0xED) packed BGR, matches the workbook's HyperLink style

function Set-HyperlinkCell {
    param($ws, [string]$cellAddr, [string]$displayText)

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $hyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
    $font = $ws.Range($cellAddr).Font
    $font.Underline = $true
    $font.Color = $hyperlinkColor
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Set-HyperlinkCell $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-02 01:08:45"

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HyperlinkCell $wsZhCn "A2" "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 01:08:39"

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HyperlinkCell $wsDeDe "A2" "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 01:08:45"
